# Sync attendance_reports, modules_schedules, and assets from main repo - 2026-01-22 06:26:48
#
# The "Recorded By" column (G) lists session recorders as a comma-joined
# string. Wherever a session was recorded by both the system auto-fill and
# the user dnasr281@gmail.com, the two names need to swap order so "System"
# is listed first: "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$lastRow = $ws.Cells(1, 1).End(4).Row
$recordedByCol = 7  # Column G - "Recorded By"

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $recordedByCol)
    if ($cell.Text -eq $oldValue) {
        $cell.Value = $newValue
    }
}
